# Fix DateTime format
#
# The "DateTime" column (C) on the "Tableless" sheet only stored whole-day
# serial dates. Two of the sample rows actually have a time-of-day
# component, so give C3/C4 a fractional serial value and switch their
# number format to the date+time format already used elsewhere in this
# workbook for the same kind of data. Column C is widened so the longer
# formatted text still fits, and the sheet/selection that was left active
# from authoring is moved onto the cells that were just fixed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tableless")

# C3: 2023-10-15 -> 2023-10-15 18:00
$ws.Range("C3").Value = 45214.75
# C4: 2023-10-16 -> 2023-10-16 20:00
$ws.Range("C4").Value = 45215.833333333336

# Apply the date+time display format to the two corrected cells.
$ws.Range("C3:C4").NumberFormat = "d/m/yy\ h:mm;@"

# Widen column C to fit the new, longer date+time text.
$ws.Columns.Item(3).ColumnWidth = 12.5

# "Tableless" becomes the active sheet, with the fixed cells selected
# (previously "WithTable_Duplicate" was the active/selected sheet).
$ws.Activate() | Out-Null
$ws.Range("C3:C4").Select() | Out-Null
